# "Natmi following Dr Hou advice"
#
# The NATMI ligand-receptor edge table (Lama1 -> Itga2) is regenerated with
# two additional cell clusters added to the analysis: "ECs" (endothelial
# cells) as a sending cluster and "M2" (macrophages) as a target cluster.
# This grows the 2 sending x 3 target (6-row) table into a 3 sending x 4
# target (12-row) table, and every numeric column is recomputed against the
# enlarged cluster set.
#
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the various expression/specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Lama1", "Itga2", "ECs",  1, 0.3333333333333333,  0.02283333333333333,  0.06850000000000001, 0.05477178157813095, 0.05477178157813096, 3, 1,                  1.740822,            5.222466,           0.4863878955914668,  0.4863878955914669,  0.039748769,               0.357738921,        0.02664033157958258,  0.02664033157958259),
    @("ECs",  "Lama1", "Itga2", "FAPs", 1, 0.3333333333333333,  0.02283333333333333,  0.06850000000000001, 0.05477178157813095, 0.05477178157813096, 3, 1,                  1.077748,            3.233244,           0.3011241710513264,  0.3011241710513265,  0.02460857933333333,      0.221477214,        0.01649310732471899,  0.016493107324719),
    @("ECs",  "Lama1", "Itga2", "M2",   1, 0.3333333333333333,  0.02283333333333333,  0.06850000000000001, 0.05477178157813095, 0.05477178157813096, 1, 0.3333333333333333, 0.03488166666666666, 0.104645,           0.009745982326006345, 0.009745982326006345, 0.0007964647222222222,    0.007168182500000001, 0.0005338048152243442, 0.0005338048152243442),
    @("ECs",  "Lama1", "Itga2", "sCs",  1, 0.3333333333333333,  0.02283333333333333,  0.06850000000000001, 0.05477178157813095, 0.05477178157813096, 3, 1,                  0.7256300000000001,  2.17689,            0.2027419510312003,  0.2027419510312003,  0.01656855166666667,      0.149116965,        0.01110453785860502,  0.01110453785860503),
    @("FAPs", "Lama1", "Itga2", "ECs",  3, 1,                   0.3395593333333333,   1.018678,            0.8145227578751427,  0.8145227578751428,  3, 1,                  1.740822,            5.222466,           0.4863878955914668,  0.4863878955914669,  0.5911123577719999,       5.320011219947999,  0.3961740101142485,   0.3961740101142486),
    @("FAPs", "Lama1", "Itga2", "FAPs", 3, 1,                   0.3395593333333333,   1.018678,            0.8145227578751427,  0.8145227578751428,  3, 1,                  1.077748,            3.233244,           0.3011241710513264,  0.3011241710513265,  0.3659593923813333,       3.293634531432,     0.2452724902675926,   0.2452724902675927),
    @("FAPs", "Lama1", "Itga2", "M2",   3, 1,                   0.3395593333333333,   1.018678,            0.8145227578751427,  0.8145227578751428,  1, 0.3333333333333333, 0.03488166666666666, 0.104645,           0.009745982326006345, 0.009745982326006345, 0.01184439547888889,      0.10659955931,      0.007938324402381087, 0.007938324402381089),
    @("FAPs", "Lama1", "Itga2", "sCs",  3, 1,                   0.3395593333333333,   1.018678,            0.8145227578751427,  0.8145227578751428,  3, 1,                  0.7256300000000001,  2.17689,            0.2027419510312003,  0.2027419510312003,  0.2463944390466667,       2.21754995142,      0.1651379330909204,   0.1651379330909204),
    @("sCs",  "Lama1", "Itga2", "ECs",  3, 1,                   0.05448866666666666,  0.163466,            0.1307054605467263,  0.1307054605467264,  3, 1,                  1.740822,            5.222466,           0.4863878955914668,  0.4863878955914669,  0.094855069684,           0.853695627156,     0.06357355389763571,  0.06357355389763573),
    @("sCs",  "Lama1", "Itga2", "FAPs", 3, 1,                   0.05448866666666666,  0.163466,            0.1307054605467263,  0.1307054605467264,  3, 1,                  1.077748,            3.233244,           0.3011241710513264,  0.3011241710513265,  0.05872505152266666,      0.528525463704,     0.03935857345901482,  0.03935857345901483),
    @("sCs",  "Lama1", "Itga2", "M2",   3, 1,                   0.05448866666666666,  0.163466,            0.1307054605467263,  0.1307054605467264,  1, 0.3333333333333333, 0.03488166666666666, 0.104645,           0.009745982326006345, 0.009745982326006345, 0.001900655507777778,     0.01710589957,      0.001273853108400914, 0.001273853108400915),
    @("sCs",  "Lama1", "Itga2", "sCs",  3, 1,                   0.05448866666666666,  0.163466,            0.1307054605467263,  0.1307054605467264,  3, 1,                  0.7256300000000001,  2.17689,            0.2027419510312003,  0.2027419510312003,  0.03953861119333334,      0.35584750074,      0.02649948008167487,  0.02649948008167488)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowValues[$j]
    }
}

$wb.Save()
